# Balanced IO card (AM RoW variant) - update de-emphasis circuit BoM rows
# to match the pre-emphasis changes made in a previous commit.
#
# Row 4  (C7 C10 - de-emphasis film capacitor): 2n7 -> 5n6
# Row 13 (R7 R15 - de-emphasis resistor):        20k -> 8k2
# Row 14 (R8 R16 - de-emphasis resistor):        7k5 -> 5k1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Value" column first for all three changed rows (R7 R15, R8 R16,
# then C7 C10), then fill in the remaining detail columns row by row -
# mirrors the order the new shared strings appear in the saved workbook.

# Row 13: R7 R15 - resistor value 20k -> 8k2
$ws.Range("C13").Value = "8k2"

# Row 14: R8 R16 - resistor value 7k5 -> 5k1
$ws.Range("C14").Value = "5k1"

# Row 4: C7 C10 - capacitor value 2n7 -> 5n6
$ws.Range("C4").Value = "5n6"
$ws.Range("E4").Value = "5.6nF 63V 5% film"
$ws.Range("G4").Value = "MMK5562J63J01L16.5TA18"
$ws.Range("I4").Value = "80-MMK5562J63J01TA18"
$ws.Range("J4").Value = 0.229

# Row 13 continued: R7 R15 detail columns
$ws.Range("E13").Value = "8.2k 0.6W 1% metal film"
$ws.Range("G13").Value = "MCMF006FF8201A50"
$ws.Range("I13").Value = "2401778"
$ws.Range("J13").Value = 0.0284

# Row 14 continued: R8 R16 detail columns
$ws.Range("E14").Value = "5.1k 0.6W 1% metal film"
$ws.Range("G14").Value = "MCMF006FF5101A50"
$ws.Range("I14").Value = "2401773"
$ws.Range("J14").Value = 0.0379
